$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10.. down by one.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new HVAC category entry.
$ws.Range("A10").Value = "HVAC"
$ws.Range("B10").Value = "Packaged Air Conditioner"
$ws.Range("C10").Value = "SKM"
$ws.Range("D10").Value = "APMR-A"

Write-Host "New row 10:" $ws.Range("A10").Text "|" $ws.Range("B10").Text "|" $ws.Range("C10").Text "|" $ws.Range("D10").Text
Write-Host "Used range now:" $ws.UsedRange.Address()

# Restore the cursor/selection like the re-saved workbook (cosmetic, matches diff).
$ws.Range("G15").Select()


